$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Cells.Item(4, 7).Value = 2.25
$ws.Cells.Item(4, 8).Value = 2.88
$ws.Cells.Item(4, 9).Value = 3.8
$ws.Cells.Item(4, 10).Value = 3.1
$ws.Cells.Item(4, 11).Value = 1.91
$ws.Cells.Item(4, 12).Value = 4.5
$ws.Cells.Item(4, 15).Value = 1.53
$ws.Cells.Item(4, 16).Value = 2.38
$ws.Cells.Item(4, 17).Value = 2.7
$ws.Cells.Item(4, 18).Value = 1.44
$ws.Cells.Item(4, 19).Value = 5.5
$ws.Cells.Item(4, 20).Value = 1.14
$ws.Cells.Item(4, 21).Value = 1.62
$ws.Cells.Item(4, 22).Value = 2.2
$ws.Cells.Item(4, 23).Value = 2.2
$ws.Cells.Item(4, 24).Value = 1.62
$ws.Cells.Item(4, 26).Value = 9
$ws.Cells.Item(4, 27).Value = 10
$ws.Cells.Item(4, 28).Value = 21
$ws.Cells.Item(4, 29).Value = 23
$ws.Cells.Item(4, 31).Value = 6
$ws.Cells.Item(4, 33).Value = 19
$ws.Cells.Item(4, 36).Value = 8
$ws.Cells.Item(4, 37).Value = 17
$ws.Cells.Item(4, 44).Value = 4.6
$ws.Cells.Item(4, 45).Value = 1.2

# Row 5
$ws.Cells.Item(5, 7).Value = 2.35
$ws.Cells.Item(5, 9).Value = 3.4
$ws.Cells.Item(5, 10).Value = 3.2
$ws.Cells.Item(5, 12).Value = 4
$ws.Cells.Item(5, 19).Value = 5.5
$ws.Cells.Item(5, 20).Value = 1.14
$ws.Cells.Item(5, 26).Value = 10
$ws.Cells.Item(5, 28).Value = 23
$ws.Cells.Item(5, 29).Value = 23
$ws.Cells.Item(5, 36).Value = 7.5
$ws.Cells.Item(5, 37).Value = 15
$ws.Cells.Item(5, 42).Value = 2
$ws.Cells.Item(5, 43).Value = 1.85

# Row 6
$ws.Cells.Item(6, 7).Value = 1.83
$ws.Cells.Item(6, 8).Value = 3.2
$ws.Cells.Item(6, 10).Value = 2.6
$ws.Cells.Item(6, 11).Value = 1.91
$ws.Cells.Item(6, 15).Value = 1.57
$ws.Cells.Item(6, 16).Value = 2.25
$ws.Cells.Item(6, 19).Value = 6
$ws.Cells.Item(6, 20).Value = 1.13
$ws.Cells.Item(6, 21).Value = 1.62
$ws.Cells.Item(6, 22).Value = 2.2
$ws.Cells.Item(6, 23).Value = 2.5
$ws.Cells.Item(6, 24).Value = 1.5
$ws.Cells.Item(6, 25).Value = 4.75
$ws.Cells.Item(6, 28).Value = 15
$ws.Cells.Item(6, 29).Value = 21
$ws.Cells.Item(6, 31).Value = 6
$ws.Cells.Item(6, 33).Value = 23
$ws.Cells.Item(6, 34).Value = 101
$ws.Cells.Item(6, 36).Value = 9.5
$ws.Cells.Item(6, 41).Value = 67
$ws.Cells.Item(6, 42).Value = 2.1
$ws.Cells.Item(6, 43).Value = 1.78
$ws.Cells.Item(6, 44).Value = 4.8
$ws.Cells.Item(6, 45).Value = 1.19

# Row 7
$ws.Cells.Item(7, 7).Value = 1.9
$ws.Cells.Item(7, 8).Value = 3.1
$ws.Cells.Item(7, 9).Value = 5
$ws.Cells.Item(7, 10).Value = 2.75
$ws.Cells.Item(7, 17).Value = 3.4
$ws.Cells.Item(7, 18).Value = 1.33
$ws.Cells.Item(7, 21).Value = 1.73
$ws.Cells.Item(7, 22).Value = 2.08
$ws.Cells.Item(7, 29).Value = 23
$ws.Cells.Item(7, 30).Value = 51
$ws.Cells.Item(7, 31).Value = 5
$ws.Cells.Item(7, 32).Value = 6.5
$ws.Cells.Item(7, 45).Value = 1.14

# Row 11
$ws.Cells.Item(11, 7).Value = 1.7
$ws.Cells.Item(11, 8).Value = 3.4
$ws.Cells.Item(11, 19).Value = 4.5
$ws.Cells.Item(11, 20).Value = 1.18
$ws.Cells.Item(11, 23).Value = 2.38
$ws.Cells.Item(11, 24).Value = 1.53
$ws.Cells.Item(11, 25).Value = 5
$ws.Cells.Item(11, 27).Value = 9.5
$ws.Cells.Item(11, 28).Value = 13
$ws.Cells.Item(11, 31).Value = 6.5
$ws.Cells.Item(11, 33).Value = 23
$ws.Cells.Item(11, 36).Value = 10
$ws.Cells.Item(11, 42).Value = 1.85
$ws.Cells.Item(11, 43).Value = 2
$ws.Cells.Item(11, 44).Value = 3.8
$ws.Cells.Item(11, 45).Value = 1.25

# Row 29
$ws.Cells.Item(29, 7).Value = 1.48
$ws.Cells.Item(29, 8).Value = 4.5
$ws.Cells.Item(29, 9).Value = 6.5
$ws.Cells.Item(29, 10).Value = 2
$ws.Cells.Item(29, 11).Value = 2.38
$ws.Cells.Item(29, 12).Value = 6
$ws.Cells.Item(29, 13).Value = 1.04
$ws.Cells.Item(29, 14).Value = 13
$ws.Cells.Item(29, 15).Value = 1.22
$ws.Cells.Item(29, 16).Value = 4
$ws.Cells.Item(29, 17).Value = 1.73
$ws.Cells.Item(29, 18).Value = 2.08
$ws.Cells.Item(29, 19).Value = 2.75
$ws.Cells.Item(29, 20).Value = 1.4
$ws.Cells.Item(29, 21).Value = 1.33
$ws.Cells.Item(29, 22).Value = 3.25

# Row 30
$ws.Cells.Item(30, 7).Value = 2.35
$ws.Cells.Item(30, 8).Value = 3.3
$ws.Cells.Item(30, 17).Value = 2.1
$ws.Cells.Item(30, 18).Value = 1.7
$ws.Cells.Item(30, 28).Value = 21
$ws.Cells.Item(30, 32).Value = 6.5

# Row 31
$ws.Cells.Item(31, 9).Value = 2.75
$ws.Cells.Item(31, 12).Value = 3.5
$ws.Cells.Item(31, 25).Value = 7.5
$ws.Cells.Item(31, 37).Value = 12

# Row 32
$ws.Cells.Item(32, 7).Value = 2.88
$ws.Cells.Item(32, 9).Value = 2.55
$ws.Cells.Item(32, 10).Value = 3.75
$ws.Cells.Item(32, 11).Value = 1.95
$ws.Cells.Item(32, 12).Value = 3.4
$ws.Cells.Item(32, 17).Value = 2.4
$ws.Cells.Item(32, 18).Value = 1.53
$ws.Cells.Item(32, 19).Value = 4.5
$ws.Cells.Item(32, 20).Value = 1.18
$ws.Cells.Item(32, 21).Value = 1.53
$ws.Cells.Item(32, 22).Value = 2.38
$ws.Cells.Item(32, 23).Value = 2
$ws.Cells.Item(32, 24).Value = 1.73
$ws.Cells.Item(32, 42).Value = 1.83
$ws.Cells.Item(32, 43).Value = 1.98
